$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-23 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-24 Friday", 2) | Out-Null
$d.Content.Find.Execute("68-14=54", $true, $false, $false, $false, $false, $true, 1, $false, "17+79=96", 2) | Out-Null
$d.Content.Find.Execute("86-9=77", $true, $false, $false, $false, $false, $true, 1, $false, "5+5=10", 2) | Out-Null
$d.Content.Find.Execute("62-48=14", $true, $false, $false, $false, $false, $true, 1, $false, "86-51=35", 2) | Out-Null
$d.Content.Find.Execute("97-25=72", $true, $false, $false, $false, $false, $true, 1, $false, "59-32=27", 2) | Out-Null
$d.Content.Find.Execute("25+10=35", $true, $false, $false, $false, $false, $true, 1, $false, "98-87=11", 2) | Out-Null
$d.Content.Find.Execute("86-20=66", $true, $false, $false, $false, $false, $true, 1, $false, "71+4=75", 2) | Out-Null
$d.Content.Find.Execute("6+44=50", $true, $false, $false, $false, $false, $true, 1, $false, "38-6=32", 2) | Out-Null
$d.Content.Find.Execute("18+17=35", $true, $false, $false, $false, $false, $true, 1, $false, "50+14=64", 2) | Out-Null
$d.Content.Find.Execute("34-16=18", $true, $false, $false, $false, $false, $true, 1, $false, "88-3=85", 2) | Out-Null
$d.Content.Find.Execute("49-10=39", $true, $false, $false, $false, $false, $true, 1, $false, "80-34=46", 2) | Out-Null
$d.Content.Find.Execute("55-42=13", $true, $false, $false, $false, $false, $true, 1, $false, "98-30=68", 2) | Out-Null
$d.Content.Find.Execute("67-44=23", $true, $false, $false, $false, $false, $true, 1, $false, "15+39=54", 2) | Out-Null
$d.Content.Find.Execute("63-7=56", $true, $false, $false, $false, $false, $true, 1, $false, "3+92=95", 2) | Out-Null
$d.Content.Find.Execute("71+20=91", $true, $false, $false, $false, $false, $true, 1, $false, "94-72=22", 2) | Out-Null
$d.Content.Find.Execute("51-39=12", $true, $false, $false, $false, $false, $true, 1, $false, "51+5=56", 2) | Out-Null
$d.Content.Find.Execute("12+48=60", $true, $false, $false, $false, $false, $true, 1, $false, "43+7=50", 2) | Out-Null
$d.Content.Find.Execute("27-15=12", $true, $false, $false, $false, $false, $true, 1, $false, "64+30=94", 2) | Out-Null
$d.Content.Find.Execute("57+8=65", $true, $false, $false, $false, $false, $true, 1, $false, "49-1=48", 2) | Out-Null
$d.Content.Find.Execute("85-41=44", $true, $false, $false, $false, $false, $true, 1, $false, "27+0=27", 2) | Out-Null
$d.Content.Find.Execute("79-41=38", $true, $false, $false, $false, $false, $true, 1, $false, "96-27=69", 2) | Out-Null
$d.Content.Find.Execute("1+44=45", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=18", 2) | Out-Null
$d.Content.Find.Execute("58-13=45", $true, $false, $false, $false, $false, $true, 1, $false, "52+14=66", 2) | Out-Null
$d.Content.Find.Execute("15+72=87", $true, $false, $false, $false, $false, $true, 1, $false, "62-16=46", 2) | Out-Null
$d.Content.Find.Execute("48+30=78", $true, $false, $false, $false, $false, $true, 1, $false, "92-74=18", 2) | Out-Null
$d.Content.Find.Execute("41-3=38", $true, $false, $false, $false, $false, $true, 1, $false, "46-18=28", 2) | Out-Null
$d.Content.Find.Execute("66-50=16", $true, $false, $false, $false, $false, $true, 1, $false, "54+35=89", 2) | Out-Null
$d.Content.Find.Execute("88-71=17", $true, $false, $false, $false, $false, $true, 1, $false, "99-52=47", 2) | Out-Null
$d.Content.Find.Execute("21-4=17", $true, $false, $false, $false, $false, $true, 1, $false, "96-6=90", 2) | Out-Null
$d.Content.Find.Execute("46+45=91", $true, $false, $false, $false, $false, $true, 1, $false, "36+36=72", 2) | Out-Null
$d.Content.Find.Execute("48-44=4", $true, $false, $false, $false, $false, $true, 1, $false, "68-62=6", 2) | Out-Null
$d.Content.Find.Execute("12+37=49", $true, $false, $false, $false, $false, $true, 1, $false, "23-1=22", 2) | Out-Null
$d.Content.Find.Execute("6-6=0", $true, $false, $false, $false, $false, $true, 1, $false, "91+0=91", 2) | Out-Null
$d.Content.Find.Execute("16+70=86", $true, $false, $false, $false, $false, $true, 1, $false, "97-49=48", 2) | Out-Null
$d.Content.Find.Execute("66+22=88", $true, $false, $false, $false, $false, $true, 1, $false, "33+62=95", 2) | Out-Null
$d.Content.Find.Execute("67+26=93", $true, $false, $false, $false, $false, $true, 1, $false, "12+50=62", 2) | Out-Null
$d.Content.Find.Execute("22+76=98", $true, $false, $false, $false, $false, $true, 1, $false, "14+1=15", 2) | Out-Null
$d.Content.Find.Execute("82-55=27", $true, $false, $false, $false, $false, $true, 1, $false, "57-36=21", 2) | Out-Null
$d.Content.Find.Execute("24-11=13", $true, $false, $false, $false, $false, $true, 1, $false, "42-15=27", 2) | Out-Null
$d.Content.Find.Execute("56+4=60", $true, $false, $false, $false, $false, $true, 1, $false, "55+16=71", 2) | Out-Null
$d.Content.Find.Execute("79-52=27", $true, $false, $false, $false, $false, $true, 1, $false, "50-16=34", 2) | Out-Null
$d.Content.Find.Execute("38+26=64", $true, $false, $false, $false, $false, $true, 1, $false, "17-4=13", 2) | Out-Null
$d.Content.Find.Execute("90-25=65", $true, $false, $false, $false, $false, $true, 1, $false, "51+23=74", 2) | Out-Null
$d.Content.Find.Execute("33-21=12", $true, $false, $false, $false, $false, $true, 1, $false, "27-17=10", 2) | Out-Null
$d.Content.Find.Execute("57+12=69", $true, $false, $false, $false, $false, $true, 1, $false, "16+75=91", 2) | Out-Null
$d.Content.Find.Execute("64-12=52", $true, $false, $false, $false, $false, $true, 1, $false, "54+7=61", 2) | Out-Null
$d.Content.Find.Execute("66-32=34", $true, $false, $false, $false, $false, $true, 1, $false, "60+17=77", 2) | Out-Null
$d.Content.Find.Execute("97-41=56", $true, $false, $false, $false, $false, $true, 1, $false, "39+22=61", 2) | Out-Null
$d.Content.Find.Execute("44-8=36", $true, $false, $false, $false, $false, $true, 1, $false, "92-85=7", 2) | Out-Null
$d.Content.Find.Execute("37+33=70", $true, $false, $false, $false, $false, $true, 1, $false, "6+6=12", 2) | Out-Null
$d.Content.Find.Execute("5+75=80", $true, $false, $false, $false, $false, $true, 1, $false, "35-12=23", 2) | Out-Null
$d.Content.Find.Execute("94-79=15", $true, $false, $false, $false, $false, $true, 1, $false, "9+77=86", 2) | Out-Null
$d.Content.Find.Execute("0+78=78", $true, $false, $false, $false, $false, $true, 1, $false, "73-61=12", 2) | Out-Null
$d.Content.Find.Execute("47+36=83", $true, $false, $false, $false, $false, $true, 1, $false, "55+37=92", 2) | Out-Null
$d.Content.Find.Execute("45+40=85", $true, $false, $false, $false, $false, $true, 1, $false, "98-72=26", 2) | Out-Null
$d.Content.Find.Execute("20-16=4", $true, $false, $false, $false, $false, $true, 1, $false, "78+11=89", 2) | Out-Null
$d.Content.Find.Execute("42-31=11", $true, $false, $false, $false, $false, $true, 1, $false, "11+57=68", 2) | Out-Null
$d.Content.Find.Execute("40+11=51", $true, $false, $false, $false, $false, $true, 1, $false, "78-20=58", 2) | Out-Null
$d.Content.Find.Execute("92+5=97", $true, $false, $false, $false, $false, $true, 1, $false, "1+9=10", 2) | Out-Null
$d.Content.Find.Execute("21+41=62", $true, $false, $false, $false, $false, $true, 1, $false, "20+71=91", 2) | Out-Null
$d.Content.Find.Execute("94-11=83", $true, $false, $false, $false, $false, $true, 1, $false, "96-22=74", 2) | Out-Null
$d.Content.Find.Execute("13-2=11", $true, $false, $false, $false, $false, $true, 1, $false, "93+6=99", 2) | Out-Null
$d.Content.Find.Execute("0+49=49", $true, $false, $false, $false, $false, $true, 1, $false, "54+29=83", 2) | Out-Null
$d.Content.Find.Execute("12-12=0", $true, $false, $false, $false, $false, $true, 1, $false, "64+33=97", 2) | Out-Null
$d.Content.Find.Execute("86-26=60", $true, $false, $false, $false, $false, $true, 1, $false, "10-4=6", 2) | Out-Null
$d.Content.Find.Execute("64+29=93", $true, $false, $false, $false, $false, $true, 1, $false, "11-3=8", 2) | Out-Null
$d.Content.Find.Execute("1+86=87", $true, $false, $false, $false, $false, $true, 1, $false, "68-15=53", 2) | Out-Null
$d.Content.Find.Execute("46-26=20", $true, $false, $false, $false, $false, $true, 1, $false, "12+33=45", 2) | Out-Null
$d.Content.Find.Execute("22+21=43", $true, $false, $false, $false, $false, $true, 1, $false, "40+43=83", 2) | Out-Null
$d.Content.Find.Execute("95-20=75", $true, $false, $false, $false, $false, $true, 1, $false, "18+12=30", 2) | Out-Null
$d.Content.Find.Execute("69-18=51", $true, $false, $false, $false, $false, $true, 1, $false, "68-23=45", 2) | Out-Null
$d.Content.Find.Execute("4+45=49", $true, $false, $false, $false, $false, $true, 1, $false, "64-31=33", 2) | Out-Null
$d.Content.Find.Execute("98-58=40", $true, $false, $false, $false, $false, $true, 1, $false, "59-38=21", 2) | Out-Null
$d.Content.Find.Execute("36-10=26", $true, $false, $false, $false, $false, $true, 1, $false, "76+15=91", 2) | Out-Null
$d.Content.Find.Execute("88-44=44", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=0", 2) | Out-Null
$d.Content.Find.Execute("72-48=24", $true, $false, $false, $false, $false, $true, 1, $false, "44+41=85", 2) | Out-Null
$d.Content.Find.Execute("62-51=11", $true, $false, $false, $false, $false, $true, 1, $false, "6+34=40", 2) | Out-Null
$d.Content.Find.Execute("46-37=9", $true, $false, $false, $false, $false, $true, 1, $false, "23+40=63", 2) | Out-Null
$d.Content.Find.Execute("61-20=41", $true, $false, $false, $false, $false, $true, 1, $false, "27+17=44", 2) | Out-Null
$d.Content.Find.Execute("87+0=87", $true, $false, $false, $false, $false, $true, 1, $false, "42+2=44", 2) | Out-Null
$d.Content.Find.Execute("44-19=25", $true, $false, $false, $false, $false, $true, 1, $false, "22+57=79", 2) | Out-Null
$d.Content.Find.Execute("42+8=50", $true, $false, $false, $false, $false, $true, 1, $false, "31-0=31", 2) | Out-Null
$d.Content.Find.Execute("38-29=9", $true, $false, $false, $false, $false, $true, 1, $false, "77-60=17", 2) | Out-Null
$d.Content.Find.Execute("32+52=84", $true, $false, $false, $false, $false, $true, 1, $false, "35+53=88", 2) | Out-Null
$d.Content.Find.Execute("54+18=72", $true, $false, $false, $false, $false, $true, 1, $false, "3+30=33", 2) | Out-Null
$d.Content.Find.Execute("76-52=24", $true, $false, $false, $false, $false, $true, 1, $false, "61-52=9", 2) | Out-Null
$d.Content.Find.Execute("41+32=73", $true, $false, $false, $false, $false, $true, 1, $false, "14+1=15", 2) | Out-Null
$d.Content.Find.Execute("48-38=10", $true, $false, $false, $false, $false, $true, 1, $false, "58-33=25", 2) | Out-Null
$d.Content.Find.Execute("35+52=87", $true, $false, $false, $false, $false, $true, 1, $false, "23-20=3", 2) | Out-Null
$d.Content.Find.Execute("14+69=83", $true, $false, $false, $false, $false, $true, 1, $false, "73-25=48", 2) | Out-Null
$d.Content.Find.Execute("4+77=81", $true, $false, $false, $false, $false, $true, 1, $false, "46-6=40", 2) | Out-Null
$d.Content.Find.Execute("44+3=47", $true, $false, $false, $false, $false, $true, 1, $false, "99-23=76", 2) | Out-Null
$d.Content.Find.Execute("1+53=54", $true, $false, $false, $false, $false, $true, 1, $false, "78-7=71", 2) | Out-Null
$d.Content.Find.Execute("72-3=69", $true, $false, $false, $false, $false, $true, 1, $false, "65-59=6", 2) | Out-Null
$d.Content.Find.Execute("82-33=49", $true, $false, $false, $false, $false, $true, 1, $false, "89-81=8", 2) | Out-Null
$d.Content.Find.Execute("39-14=25", $true, $false, $false, $false, $false, $true, 1, $false, "91+1=92", 2) | Out-Null
$d.Content.Find.Execute("49+37=86", $true, $false, $false, $false, $false, $true, 1, $false, "64-41=23", 2) | Out-Null
$d.Content.Find.Execute("23+51=74", $true, $false, $false, $false, $false, $true, 1, $false, "92-78=14", 2) | Out-Null
$d.Content.Find.Execute("9+62=71", $true, $false, $false, $false, $false, $true, 1, $false, "48+27=75", 2) | Out-Null
$d.Content.Find.Execute("59-54=5", $true, $false, $false, $false, $false, $true, 1, $false, "38+36=74", 2) | Out-Null
$d.Content.Find.Execute("47+7=54", $true, $false, $false, $false, $false, $true, 1, $false, "82-39=43", 2) | Out-Null

Write-Host "Replacements complete: 101"
